# Edit the "Return of Biopsy Results within 1 Month" row of the RoR
# (Biospecimens) table: widen the row, clarify the numerator wording,
# and rename the "Biopsy date" column header to "Past a month since
# biopsy" to reflect the new divisor language ("change in RoR divisor").

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$row = $t.Rows.Item(25)

# Sanity-check we are on the right row before mutating it.
$label = $row.Cells.Item(1).Range.Text
if ($label -notlike "Return of Biopsy Results within 1 Month*") {
    throw "Unexpected row; aborting (label was: $label)"
}

# 1) Row grows from 907 -> 1051 twips (twips = points * 20).
$row.Height = 1051 / 20

# 2) Append " past a month from biopsy" to the end of the numerator
#    description cell (cell 2), just before the paragraph end mark.
$cell2 = $row.Cells.Item(2)
$insertPoint = $d.Range($cell2.Range.End - 1, $cell2.Range.End - 1)
$insertPoint.InsertAfter(" past a month from biopsy")

# 3) Replace "Biopsy date" with "Past a month since biopsy" in cell 3.
$cell3 = $row.Cells.Item(3)
$cell3Text = $d.Range($cell3.Range.Start, $cell3.Range.End - 1)
$cell3Text.Text = "Past a month since biopsy"

Write-Host ("Cell2 now: " + $cell2.Range.Text)
Write-Host ("Cell3 now: " + $cell3.Range.Text)
Write-Host ("Row height now (twips): " + ($row.Height * 20))
